# Reorder the list of mode strings so that "#ignore" and "#aliases"
# sit in the middle of the list (between "#config hidden()" and
# "#meta hidden()"), swapping places with "#meta hidden()" and
# "#notes hidden()".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "#ignore"
$ws.Range("A8").Value = "#aliases"
$ws.Range("A9").Value = "#meta hidden()"
$ws.Range("A10").Value = "#notes hidden()"
